$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new source / DOI row below the existing two data rows.
# The DOI link is entered before the source name, matching the order
# the new strings were registered in the shared-string table.
$ws.Range("B3").Value = "https://doi.org/10.1016/S1146-609X(00)01084-5"
$ws.Range("A3").Value = "PérezFernández_2000"

# Leave the freshly-entered cell selected, as in the saved workbook.
$ws.Range("B3").Select() | Out-Null

# Align page setup with the other Kew datasets (A4 portrait).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
